$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 454, pushing the existing rows 454:572
# down to 455:573 (carries their values/formatting along).
$ws.Rows.Item(454).Insert()

# Populate the newly inserted row 454 with the new weekly data point.
$ws.Cells.Item(454, 1).Value = 6
$ws.Cells.Item(454, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(454, 3).Value = "Metropolitana"
$ws.Cells.Item(454, 4).Value = 44754
$ws.Cells.Item(454, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(454, 5).Value = 13
$ws.Cells.Item(454, 6).Value = 100112012
$ws.Cells.Item(454, 7).Value = "Espinaca"
$ws.Cells.Item(454, 8).Value = "Sin especificar"
$ws.Cells.Item(454, 9).Value = "Primera"
$ws.Cells.Item(454, 10).Value = 420
$ws.Cells.Item(454, 11).Value = 7500
$ws.Cells.Item(454, 12).Value = 8000
$ws.Cells.Item(454, 13).Value = 7714
$ws.Cells.Item(454, 14).Value = '$/cuna 10 kilos'
$ws.Cells.Item(454, 15).Value = "Región Metropolitana"
$ws.Cells.Item(454, 16).Value = 771
$ws.Cells.Item(454, 17).Value = 10
$ws.Cells.Item(454, 18).Value = "Hortaliza"
